$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" everywhere it appears ---
# Overview sheet: zh-cn / de-de status columns (E, F) for both data rows (2, 3)
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Range("E3").Value = "In Translation"
$ws1.Range("F3").Value = "In Translation"

# zh-cn sheet: Status column (C) for both data rows
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "In Translation"
$ws2.Range("C3").Value = "In Translation"

# de-de sheet: Status column (C) for both data rows
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "In Translation"
$ws3.Range("C3").Value = "In Translation"

# --- Narrow the status columns (report regenerated with shorter "In Translation" text) ---
# Overview: columns E and F (zh-cn / de-de status)
$ws1.Columns("E:F").ColumnWidth = 12.5

# zh-cn / de-de: column C (Status)
$ws2.Columns("C:C").ColumnWidth = 12.5
$ws3.Columns("C:C").ColumnWidth = 12.5
